# Fruta / hortaliza, semanal
# Insert a new weekly record as row 477 (pushing the existing rows 477-529 down to 478-530).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 477; everything currently at/after row 477 shifts down one row.
$ws.Rows("477:477").Insert()

# Populate the new row with the new weekly price record.
$ws.Cells.Item(477, 1).Value = 11
$ws.Cells.Item(477, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(477, 3).Value = "Bíobío"
$ws.Cells.Item(477, 4).Value = 44488
$ws.Cells.Item(477, 5).Value = 8
$ws.Cells.Item(477, 6).Value = "Fruta"
$ws.Cells.Item(477, 7).Value = 100108
$ws.Cells.Item(477, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(477, 9).Value = 100108006
$ws.Cells.Item(477, 10).Value = "Plátano"
$ws.Cells.Item(477, 11).Value = "Sin especificar"
$ws.Cells.Item(477, 12).Value = "Pintón"
$ws.Cells.Item(477, 13).Value = 430
$ws.Cells.Item(477, 14).Value = 24000
$ws.Cells.Item(477, 15).Value = 25000
$ws.Cells.Item(477, 16).Value = 24465
$ws.Cells.Item(477, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(477, 18).Value = "Ecuador"
$ws.Cells.Item(477, 19).Value = 1223
$ws.Cells.Item(477, 20).Value = 20
